$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 213, pushing existing rows 213:336 down to 214:337.
$ws.Rows.Item(213).Insert()

# Populate the newly inserted row 213 with the new weekly data point.
$ws.Cells.Item(213, 1).Value = 9
$ws.Cells.Item(213, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(213, 3).Value = "Metropolitana"
$ws.Cells.Item(213, 4).Value = 44582
$ws.Cells.Item(213, 5).Value = 13
$ws.Cells.Item(213, 6).Value = "Fruta"
$ws.Cells.Item(213, 7).Value = 100108
$ws.Cells.Item(213, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(213, 9).Value = 100108002
$ws.Cells.Item(213, 10).Value = "Mango"
$ws.Cells.Item(213, 11).Value = "Sin especificar"
$ws.Cells.Item(213, 12).Value = "Primera"
$ws.Cells.Item(213, 13).Value = 660
$ws.Cells.Item(213, 14).Value = 5000
$ws.Cells.Item(213, 15).Value = 5500
$ws.Cells.Item(213, 16).Value = 5265
$ws.Cells.Item(213, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(213, 18).Value = "Perú"
$ws.Cells.Item(213, 19).Value = 1316
$ws.Cells.Item(213, 20).Value = 4
